# Daily attendance processing - 2025-10-17 09:21:29
# Normalize the "Recorded By" column (G) so that entries of the form
# "<name>, System" become "System, <name>" for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.EndsWith(", System")) {
        $parts = $val -split ", "
        if ($parts.Count -eq 2 -and $parts[1] -eq "System") {
            $cell.Value = "System, " + $parts[0]
        }
    }
}
